$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F13").Value = 1315
$ws.Range("G13").Value = 116
$ws.Range("I13").Value = 162.61112
$ws.Range("J13").Value = 30.8961128
$ws.Range("K13").Value = 0.116
$ws.Range("M13").Value = 0.7018
$ws.Range("N13").Value = 0.133342
$ws.Range("O13").Value = 0.116
$ws.Range("Q13").Value = 8.31494844
$ws.Range("R13").Value = 1.5798402036
$ws.Range("S13").Value = -116
$ws.Range("U13").Value = -105.34192
$ws.Range("V13").Value = -20.0149648
$ws.Range("W13").Value = 66.28594843999998
$ws.Range("X13").Value = 12.5943302036
$ws.Range("Y13").Value = 78.88027864359998
